$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the simulated data values (rows 2-5, columns C,D,E,G,H,I,K,L,M)
$ws.Range("C2").Value = 0.06542581754549807449
$ws.Range("D2").Value = 0.08340285746176910597
$ws.Range("E2").Value = 0.08507460701350666454
$ws.Range("G2").Value = 0.50530815326913325336
$ws.Range("H2").Value = 0.53549954206272509794
$ws.Range("I2").Value = 0.46682690012958194536
$ws.Range("K2").Value = 0.68396335412022379874
$ws.Range("L2").Value = 0.68887933498988096481
$ws.Range("M2").Value = 0.68017530472697529387

$ws.Range("C3").Value = 0.21795111073291167925
$ws.Range("D3").Value = 0.2657716004170566082
$ws.Range("E3").Value = 0.26049996470420211381
$ws.Range("G3").Value = 0.66667698199173264229
$ws.Range("H3").Value = 0.66104406028723294853
$ws.Range("I3").Value = 0.67208950937490674615
$ws.Range("K3").Value = 0.6165022570187868034
$ws.Range("L3").Value = 0.71095215165482383846
$ws.Range("M3").Value = 0.714746576178336257

$ws.Range("C4").Value = 0.58468356609556870751
$ws.Range("D4").Value = 0.40250378461243119244
$ws.Range("E4").Value = 0.43174604912019598491
$ws.Range("G4").Value = 0.78416489251241527914
$ws.Range("H4").Value = 0.80785471333946323735
$ws.Range("I4").Value = 0.83870556540950769175
$ws.Range("K4").Value = 1.13826233041653690492
$ws.Range("L4").Value = 1.26672775095610612439
$ws.Range("M4").Value = 1.02540003563159709543

$ws.Range("C5").Value = 0.86282502387673076782
$ws.Range("D5").Value = 0.77762469642360887434
$ws.Range("E5").Value = 0.8328472754618321261
$ws.Range("G5").Value = 0.48335264027604857295
$ws.Range("H5").Value = 0.78493709588235882624
$ws.Range("I5").Value = 0.85109763632603163686
$ws.Range("K5").Value = 1.28498605292825995861
$ws.Range("L5").Value = 1.22586685847105592018
$ws.Range("M5").Value = 1.23476717769249466627

# Update the active selection to match the author's last selection in Excel
$ws.Range("K7:M10").Select()
